$d = $word.ActiveDocument
$enDash = [char]0x2013

# ---------------------------------------------------------------------
# Substantive content changes
# ---------------------------------------------------------------------

# 1) Glaukos end date: "Oct 2021 - Mar 2023 (Current)" -> "... Apr 2023 ..."
$d.Content.Find.Execute("Oct 2021 " + $enDash + " Mar", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Oct 2021 " + $enDash + " Apr", 2)

# 2) Databases line: "SQL, NoSQL, Mongo, Postgres" -> "SQL, Postgres, SQLite, Sequelize"
$d.Content.Find.Execute("SQL, NoSQL, Mongo, Postgres", $true, $false, $false, $false, $false,
                         $true, 1, $false, "SQL, Postgres, SQLite, Sequelize", 2)

# ---------------------------------------------------------------------
# Cosmetic cleanup: these re-type the same visible text across runs that
# were split apart only for w:proofErr spell/grammar markers, so the
# replace is a same-for-same Find/Replace that forces the editor to
# rebuild the run (dropping the now-orphaned proofErr bookkeeping), the
# same way accepting/clearing a proofing squiggle would in real Word.
# ---------------------------------------------------------------------

$d.Content.Find.Execute("github.com/sammcgrail ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "github.com/sammcgrail ", 2)

$d.Content.Find.Execute("linkedin.com/in/sammcgrail ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "linkedin.com/in/sammcgrail ", 2)

$d.Content.Find.Execute("Glaukos ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Glaukos ", 2)

$d.Content.Find.Execute("doblePRIME  ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "doblePRIME  ", 2)

$text = "Extensive use of various javascript data visualization libraries (some altered for custom business requests) to allow large datasets to be presented more clearly. D3.js based viz with SVG and HTML5, using canvas elements and WebGL to render thousands of datapoints in the browser.  "
$d.Content.Find.Execute($text, $true, $false, $false, $false, $false,
                         $true, 1, $false, $text, 2)

$d.Content.Find.Execute("Launch Academy  ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Launch Academy  ", 2)

$d.Content.Find.Execute("Junior Full Stack Web Developer  ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Junior Full Stack Web Developer  ", 2)

$d.Content.Find.Execute("Professional Musician", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Professional Musician", 2)

$d.Content.Find.Execute("Amherst, MA ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Amherst, MA ", 2)

$d.Content.Find.Execute("Node, *.js, shell-fu, git, bash, custom VM curation", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Node, *.js, shell-fu, git, bash, custom VM curation", 2)

$d.Save()
